{"js": "// Bump the font size of the two \"datos personales\" / legal-notice\n// paragraphs at the end of the document's first section from 5pt (sz=10)\n// to 5.5pt (sz=11), and normalize that section's page margins to\n// 1 inch (1440 twips / 72 points) on every side.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the two target paragraphs by their distinctive leading text,\n// rather than relying on a fixed index.\nconst items = paragraphs.items;\nlet noticeParagraphs = [];\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text || \"\";\n  if (\n    t.indexOf(\"Este documento puede contener datos personales\") !== -1 ||\n    t.indexOf(\"Atendiendo a lo establecido por el art\u00edculo 72\") !== -1\n  ) {\n    noticeParagraphs.push(items[i]);\n  }\n}\n\nfor (const p of noticeParagraphs) {\n  // w:sz / w:szCs are stored in half-points; 11 half-points = 5.5pt.\n  p.font.size = 5.5;\n  p.font.sizeBidirectional = 5.5;\n}\nawait context.sync();\n\n// The second of those two paragraphs carries the embedded sectPr that\n// closes out the document's first (non-continuous) section \u2014 that is\n// Sections.items[0] here. Update its margins to 1440 twips (72pt) on\n// every side.\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nconst firstSection = sections.items[0];\nconst pageSetup = firstSection.pageSetup;\npageSetup.topMargin = 72;\npageSetup.bottomMargin = 72;\npageSetup.leftMargin = 72;\npageSetup.rightMargin = 72;\nawait context.sync();\n", "ps1": "# Bump the font size of the two \"datos personales\" / legal-notice\n# paragraphs at the end of the document's first section from 5pt (sz=10)\n# to 5.5pt (sz=11), and normalize that section's page margins to\n# 1 inch (1440 twips / 72 points) on every side.\n\n$d = $word.ActiveDocument\n\n$target1 = $null\n$target2 = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text\n    if ($t -like \"*Este documento puede contener datos personales*\") {\n        $target1 = $p\n    }\n    if ($t -like \"*Atendiendo a lo establecido por el*\") {\n        $target2 = $p\n    }\n}\n\nif ($target1 -eq $null) { throw \"Could not locate the 'Este documento puede contener datos personales' paragraph\" }\nif ($target2 -eq $null) { throw \"Could not locate the 'Atendiendo a lo establecido por el articulo 72' paragraph\" }\n\n# w:sz / w:szCs are stored in half-points; 11 half-points = 5.5pt.\n$target1.Range.Font.Size = 5.5\n$target1.Range.Font.SizeBi = 5.5\n$target2.Range.Font.Size = 5.5\n$target2.Range.Font.SizeBi = 5.5\n\n# The second of those two paragraphs carries the embedded sectPr that\n# closes out the document's first (non-continuous) section. Update its\n# margins to 1440 twips (72pt / 1 inch) on every side.\n$sec = $target2.Range.Sections.Item(1)\n$sec.PageSetup.TopMargin = 72\n$sec.PageSetup.BottomMargin = 72\n$sec.PageSetup.LeftMargin = 72\n$sec.PageSetup.RightMargin = 72\n"}
